$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Play my 'Driving at night' playlist`"",
    "`"Play country music`"",
    "`"Play Frank Sinatra`"",
    "`"Play Planet Money on NPR One.`"",
    "`"Play Adele Radio on iHeartRadio.`"",
    "`"Navigate to Union Square, San Francisco.`"",
    "`"Directions to Philz Coffee.`"",
    "`"Drive to 1600 Amphitheatre Parkway, Mountain View.`"",
    "`"Message Cody White on Hangouts.`"",
    "`"Message Cody White on Whatsapp.`""
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

[void]$ws.Range("A18").Select()
